# Update the division problems in the table on the page.
# Each data row (1, 5, 9, 13, 17) has 5 columns of "a÷b=" expressions
# that need to be replaced with new values, matching cell position.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$replacements = @(
    @{Row=1;  Col=1; Text="89÷5="},
    @{Row=1;  Col=2; Text="62÷6="},
    @{Row=1;  Col=3; Text="55÷3="},
    @{Row=1;  Col=4; Text="60÷2="},
    @{Row=1;  Col=5; Text="19÷4="},

    @{Row=5;  Col=1; Text="91÷4="},
    @{Row=5;  Col=2; Text="74÷3="},
    @{Row=5;  Col=3; Text="30÷7="},
    @{Row=5;  Col=4; Text="56÷8="},
    @{Row=5;  Col=5; Text="42÷9="},

    @{Row=9;  Col=1; Text="55÷2="},
    @{Row=9;  Col=2; Text="87÷4="},
    @{Row=9;  Col=3; Text="63÷6="},
    @{Row=9;  Col=4; Text="18÷5="},
    @{Row=9;  Col=5; Text="31÷3="},

    @{Row=13; Col=1; Text="31÷8="},
    @{Row=13; Col=2; Text="76÷9="},
    @{Row=13; Col=3; Text="69÷6="},
    @{Row=13; Col=4; Text="96÷9="},
    @{Row=13; Col=5; Text="83÷9="},

    @{Row=17; Col=1; Text="61÷4="},
    @{Row=17; Col=2; Text="35÷9="},
    @{Row=17; Col=3; Text="24÷9="},
    @{Row=17; Col=4; Text="89÷9="},
    @{Row=17; Col=5; Text="53÷6="}
)

foreach ($r in $replacements) {
    $cell = $t.Cell($r.Row, $r.Col)
    $cell.Range.Text = $r.Text
}
